$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New header cells (Wins / Losses / Ties), styled like the rest of row 1
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the team record values for every data row (2-49)
for ($r = 2; $r -le 49; $r++) {
    $ws.Cells.Item($r, 30).Value = 90  # AD - Wins
    $ws.Cells.Item($r, 31).Value = 72  # AE - Losses
    $ws.Cells.Item($r, 32).Value = 0   # AF - Ties
}
